$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1.xml)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 13635
$ws1.Range("F12").Value = 5
$ws1.Range("F14").Value = 13628
$ws1.Range("F16").Value = 608
$ws1.Range("F17").Value = 9000
$ws1.Range("F19").Value = 8106
$ws1.Range("F24").Value = 154
$ws1.Range("F30").Value = 397
$ws1.Range("F32").Value = 197

# Sheet "全部类型" (sheet4.xml)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 13635
$ws4.Range("F12").Value = 5
$ws4.Range("F14").Value = 13628
$ws4.Range("F16").Value = 608
$ws4.Range("F17").Value = 9000
$ws4.Range("F19").Value = 8106
$ws4.Range("F24").Value = 154
$ws4.Range("F32").Value = 397
$ws4.Range("F34").Value = 197
